$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below mirrors one <c> element changed in the target OOXML diff:
# plain numeric-looking strings (e.g. "226.36") must land back in the sheet
# as literal text (matching the original inlineStr cells), not get
# auto-coerced into Excel numbers - so we briefly force text format, assign,
# then restore the default style to avoid leaving stray formatting behind.
$updates = @(
    @{Cell='D2'; Value='27.877.04'},
    @{Cell='E2'; Value='  -4.56%  '},
    @{Cell='D3'; Value='1.737.00'},
    @{Cell='E3'; Value='  -4.86%  '},
    @{Cell='E4'; Value='  -0.16%  '},
    @{Cell='D5'; Value='226.36'},
    @{Cell='E5'; Value='  -3.57%  '},
    @{Cell='D6'; Value='0.5763'},
    @{Cell='E6'; Value='  -3.80%  '},
    @{Cell='E7'; Value='  -0.18%  '},
    @{Cell='D8'; Value='0.2733'},
    @{Cell='E8'; Value='  -0.78%  '},
    @{Cell='D9'; Value='23.05'},
    @{Cell='E9'; Value='  -1.23%  '},
    @{Cell='D10'; Value='0.06617'},
    @{Cell='E10'; Value='  -4.20%  '},
    @{Cell='D11'; Value='0.07536'},
    @{Cell='E11'; Value='  -0.90%  '},
    @{Cell='D12'; Value='1.742.08'},
    @{Cell='E12'; Value='  -4.74%  '},
    @{Cell='D13'; Value='4.702'},
    @{Cell='E13'; Value='  -0.10%  '},
    @{Cell='D14'; Value='0.6013'},
    @{Cell='E14'; Value='  -3.46%  '},
    @{Cell='D15'; Value='1.973.48'},
    @{Cell='E15'; Value='  -4.82%  '},
    @{Cell='D16'; Value='74.41'},
    @{Cell='E16'; Value='  -3.27%  '},
    @{Cell='D17'; Value='0.000008668'},
    @{Cell='E17'; Value='  -10.66%  '},
    @{Cell='D18'; Value='27.859.86'},
    @{Cell='E18'; Value='  -3.90%  '},
    @{Cell='D19'; Value='5.309'},
    @{Cell='E19'; Value='  -3.85%  '},
    @{Cell='E20'; Value='  -0.21%  '},
    @{Cell='D21'; Value='204.88'},
    @{Cell='E21'; Value='  -4.53%  '},
    @{Cell='D22'; Value='11.26'},
    @{Cell='E22'; Value='  -2.28%  '},
    @{Cell='D23'; Value='6.614'},
    @{Cell='E23'; Value='  -2.86%  '},
    @{Cell='E24'; Value='  -0.22%  '},
    @{Cell='D25'; Value='150.02'},
    @{Cell='E25'; Value='  -3.72%  '},
    @{Cell='D26'; Value='8.009'},
    @{Cell='E26'; Value='  +0.90%  '},
    @{Cell='D27'; Value='0.1231'},
    @{Cell='E27'; Value='  -4.15%  '},
    @{Cell='D28'; Value='16.17'},
    @{Cell='E28'; Value='  -1.67%  '},
    @{Cell='B29'; Value='Toncoin'},
    @{Cell='C29'; Value='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'},
    @{Cell='D29'; Value='1.384'},
    @{Cell='E29'; Value='  -2.92%  '},
    @{Cell='B30'; Value='Hedera'},
    @{Cell='C30'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D30'; Value='0.06177'},
    @{Cell='E30'; Value='  -4.88%  '},
    @{Cell='D31'; Value='1.391'},
    @{Cell='E31'; Value='  -3.45%  '},
    @{Cell='D32'; Value='3.738'},
    @{Cell='E32'; Value='  -1.25%  '},
    @{Cell='D33'; Value='3.725'},
    @{Cell='E33'; Value='  -1.15%  '},
    @{Cell='D34'; Value='1.674'},
    @{Cell='E34'; Value='  -2.30%  '},
    @{Cell='D35'; Value='1.033'},
    @{Cell='E35'; Value='  -4.81%  '},
    @{Cell='D36'; Value='0.6439'},
    @{Cell='E36'; Value='  +0.59%  '},
    @{Cell='D37'; Value='2.438'},
    @{Cell='E37'; Value='  -4.16%  '},
    @{Cell='D38'; Value='2.708'},
    @{Cell='E38'; Value='  -2.14%  '},
    @{Cell='D39'; Value='0.01668'},
    @{Cell='D40'; Value='1.119.05'},
    @{Cell='E40'; Value='  -1.27%  '},
    @{Cell='D41'; Value='6.166'},
    @{Cell='E41'; Value='  -6.11%  '},
    @{Cell='D42'; Value='0.8718'},
    @{Cell='E42'; Value='  -1.70%  '},
    @{Cell='D43'; Value='1.005'},
    @{Cell='E43'; Value='  +0.05%  '},
    @{Cell='D44'; Value='100.02'},
    @{Cell='E44'; Value='  -0.62%  '},
    @{Cell='D45'; Value='1.884.79'},
    @{Cell='E45'; Value='  -5.01%  '},
    @{Cell='D46'; Value='59.28'},
    @{Cell='E46'; Value='  -4.03%  '},
    @{Cell='D47'; Value='0.00000000108'},
    @{Cell='E47'; Value='  -4.24%  '},
    @{Cell='D48'; Value='1.571'},
    @{Cell='E48'; Value='  -2.03%  '},
    @{Cell='D49'; Value='8.252'},
    @{Cell='E49'; Value='  -2.08%  '},
    @{Cell='D50'; Value='0.05374'},
    @{Cell='E50'; Value='  -2.35%  '},
    @{Cell='E51'; Value='  -2.63%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $val = $u.Value
    if ($val -match '^\s*[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?\s*$') {
        # Numeric-looking text (e.g. "226.36") - force text format first so
        # Excel doesn't silently convert it to a Number cell, then drop the
        # temporary formatting back to the default style.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
